$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 38
$data = New-Object 'object[,]' $rowCount,6

$data[0,0] = "Editorial: LAUSD students’ scores are bouncing back — thanks to teachers"
$data[0,1] = "The latest test scores for L.A. students are encouraging — but not yet time for a full-throated cheer because they are still below pre-pandemic levels."
$data[0,2] = "July 28, 2024"
$data[0,3] = "a594b684-67c7-4296-beb7-90175a6e0af1"
$data[0,4] = 0
$data[0,5] = $false

$data[1,0] = "Tim Vanderhook"
$data[1,1] = "Tim Vanderhook is a visionary entrepreneur and the driving force behind Viant Technology Inc (NASDAQ: DSP), a leading advertising technology company he co-founded with his brother Chris."
$data[1,2] = "July 28, 2024"
$data[1,3] = "8f35108c-afcf-423a-ad1e-c236304d491a"
$data[1,4] = 0
$data[1,5] = $false

$data[2,0] = "Sasha Strauss"
$data[2,1] = "Sasha Strauss is a distinguished brand strategy expert, educator and speaker renowned for his 25-year career shaping brands globally."
$data[2,2] = "July 28, 2024"
$data[2,3] = "a0a3ce09-3328-4a8d-a484-ec40f182240a"
$data[2,4] = 0
$data[2,5] = $false

$data[3,0] = "Sona Shah"
$data[3,1] = "Sona Shah founded My Private Professor (MPP) to democratize access to academic resources, leveraging her expertise in business litigation and passion for education."
$data[3,2] = "July 28, 2024"
$data[3,3] = "f5e63fc1-6457-432f-bc2c-0ad9c2a19651"
$data[3,4] = 1
$data[3,5] = $false

$data[4,0] = "Lucy Santana"
$data[4,1] = "Lucy Santana has led Girls Inc. of Orange County for over 20 years, transforming it into a leading chapter nationally."
$data[4,2] = "July 28, 2024"
$data[4,3] = "42b52a4b-ff6c-4f21-961f-6b53084af30d"
$data[4,4] = 0
$data[4,5] = $false

$data[5,0] = "Michael J. Beals"
$data[5,1] = "Dr. Michael J. Beals has served as president of Vanguard University (VU) for 10 years, guiding Orange County’s oldest four-year university to significant growth and achievement."
$data[5,2] = "July 28, 2024"
$data[5,3] = "238e78cd-3ddb-4faf-adda-dc988d4b96ce"
$data[5,4] = 0
$data[5,5] = $false

$data[6,0] = "Emil Davtyan"
$data[6,1] = "Emil Davtyan, the founder and CEO of D.Law, deserves recognition for his exceptional leadership and dedication to advocating for California workers’ rights."
$data[6,2] = "July 28, 2024"
$data[6,3] = "6ce63d7d-3b1c-41b7-95f1-486c79c28426"
$data[6,4] = 0
$data[6,5] = $false

$data[7,0] = "Prioritizing Tech Innovation & Talent Retention"
$data[7,1] = "Moderator: Brian Hegarty Principal, L.A."
$data[7,2] = "July 28, 2024"
$data[7,3] = "e69e26db-b569-443d-ae18-201df0e5d4ef"
$data[7,4] = 0
$data[7,5] = $false

$data[8,0] = "Column: 99 years after the Scopes ‘monkey trial,’ religious fundamentalism still infects our schools"
$data[8,1] = "The Scopes trial made anti-evolutionists look ridiculous, but they haven’t gone away."
$data[8,2] = "July 26, 2024"
$data[8,3] = "3e03fc44-a137-4353-a803-e775a5b8d6a3"
$data[8,4] = 0
$data[8,5] = $false

$data[9,0] = "Homeless encampment cleanups do little to change numbers of people on the street, study finds"
$data[9,1] = "Three new studies shed new light on homelessness in Los Angeles, painting a mostly grim picture, with some glimmers of hope."
$data[9,2] = "July 24, 2024"
$data[9,3] = "42ce960d-b14c-4eb3-8e11-71996ab67252"
$data[9,4] = 0
$data[9,5] = $false

$data[10,0] = "LAUSD test scores rise in math and English, positive marks after pandemic setbacks"
$data[10,1] = "Supt. Alberto Carvalho said the nation’s second-largest school system has seen across-the-board improvement in math and English scores in every grade."
$data[10,2] = "July 23, 2024"
$data[10,3] = "ecf17e92-5466-441d-93c6-75004db30934"
$data[10,4] = 0
$data[10,5] = $false

$data[11,0] = "Boiling Point: Living in Death Valley"
$data[11,1] = "Death Valley residents talk about how they deal with extreme heat every day in the summer"
$data[11,2] = "July 23, 2024"
$data[11,3] = "4b55a9fa-6af2-420a-b18c-0e7239eaf453"
$data[11,4] = 0
$data[11,5] = $false

$data[12,0] = "The 22 best spots to nerd out in L.A."
$data[12,1] = "Whether you’re fiercely into tabletop games, robots, pinball, dinosaurs, space, comic books or close-up magic, here’s where you can find your people."
$data[12,2] = "July 23, 2024"
$data[12,3] = "77c0375a-7cf1-4ce9-8cfb-7bd04e4cdbb5"
$data[12,4] = 0
$data[12,5] = $false

$data[13,0] = "Carvalho, who unplugged school AI chatbot, wants task force to tell him what went wrong"
$data[13,1] = "Independent experts will be asked to look at what went amiss with LAUSD’s AI effort and helped plan next steps in the ongoing but stalled strategy."
$data[13,2] = "July 23, 2024"
$data[13,3] = "1fe90465-e048-4778-a395-5d550b966ab4"
$data[13,4] = 0
$data[13,5] = $false

$data[14,0] = "Could AI robots with lasers make herbicides — and farm workers — obsolete?"
$data[14,1] = "A shift from harmful herbicides to intelligent robots would have far-reaching consequences for California’s `$50-billion agriculture industry."
$data[14,2] = "July 22, 2024"
$data[14,3] = "e06936f8-83c8-4160-a23c-f3db6ee37fc5"
$data[14,4] = 0
$data[14,5] = $false

$data[15,0] = "Robert T. Braithwaite, Channing Hamlet and Daniel A. Platt Share Insights on the Healthcare, Biotech and Life Sciences Landscape in 2024"
$data[15,1] = "The Healthcare, Biotech & Life Sciences Roundtable panel is produced by the L.A."
$data[15,2] = "July 22, 2024"
$data[15,3] = "d1f0fed9-a126-4f51-8882-3109149ccf5a"
$data[15,4] = 0
$data[15,5] = $false

$data[16,0] = "Usha Vance: From San Francisco corporate lawyer to MAGA’s potential second lady"
$data[16,1] = "The California native and daughter of Indian academics has devoted her life to amassing the kind of elite professional credentials that draw scorn from GOP populists."
$data[16,2] = "July 19, 2024"
$data[16,3] = "7ca088c4-2128-4f44-b0bf-2473d99473ea"
$data[16,4] = 0
$data[16,5] = $false

$data[17,0] = "These California counties endured the nation’s longest streaks of excessive heat"
$data[17,1] = "Several California communities near the Arizona border had the longest streak in the nation of days that hit 90 degrees or higher in 2022, according to new U.S. Census data."
$data[17,2] = "July 19, 2024"
$data[17,3] = "8cdc3a04-c97d-4940-9bed-c2bb0ba9705e"
$data[17,4] = 0
$data[17,5] = $false

$data[18,0] = "Review: ‘Tartuffe: Born Again’ in the American South kicks up its heels at Theatricum Botanicum"
$data[18,1] = "Freyda Thomas adapts Moliere’s ‘Tartuffe’ with an American twist in a winning production at Topanga’s beloved Will Geer Theatricum Botanicum."
$data[18,2] = "July 19, 2024"
$data[18,3] = "d313955b-74fc-49cd-952f-4984d19721ad"
$data[18,4] = 0
$data[18,5] = $false

$data[19,0] = "UC regents ban views on Israel, other political opinion from university homepages"
$data[19,1] = "UC regents voted Wednesday to ban views on political opinions from main campus homepages. Opinions may be posted on other pages but it must be made clear they aren’t official UC views."
$data[19,2] = "July 18, 2024"
$data[19,3] = "e3543ab1-6a22-4a20-b75c-f397a8b47986"
$data[19,4] = 0
$data[19,5] = $false

$data[20,0] = "For China, Trump rally shooting is more evidence of America’s demise"
$data[20,1] = "America’s biggest rival says Trump shooting symbolizes dwindling power and hypocrisy of American democracy and global leadership."
$data[20,2] = "July 17, 2024"
$data[20,3] = "f35fac03-33c7-4f85-8f51-884278b4402c"
$data[20,4] = 0
$data[20,5] = $false

$data[21,0] = "Former Stanford dean, now local council member, apologizes for affair with student"
$data[21,1] = "Palo Alto City Council Member Julie Lythcott-Haims apologized for the romantic relationship she had with a Stanford undergrad more than a decade ago when she served as university dean."
$data[21,2] = "July 15, 2024"
$data[21,3] = "de7a1752-e71d-454f-a922-692414ca5b36"
$data[21,4] = 0
$data[21,5] = $false

$data[22,0] = "Opinion: J.D. Vance’s book ‘Hillbilly Elegy’ was a con job. Don’t let it slide"
$data[22,1] = "The running mate for Donald Trump showed in his 2016 memoir that he doesn’t understand or respect the working class. The media should call him out on that."
$data[22,2] = "July 15, 2024"
$data[22,3] = "1fc1abb6-6e0a-480b-ac78-9a1fb47f161c"
$data[22,4] = 0
$data[22,5] = $false

$data[23,0] = "UC regents: Protests yes, encampments no. Campus rules must be consistently enforced"
$data[23,1] = "Rich Leib, outgoing chair of the UC Board of Regents, says encampments should be banned, but protests that follow campus rules are welcomed as free speech. Many regents, senior leaders agree."
$data[23,2] = "July 15, 2024"
$data[23,3] = "17bc38cc-cf69-4560-a1c1-b31329345e0d"
$data[23,4] = 0
$data[23,5] = $false

$data[24,0] = "There’s no crying in baseball, but there are plenty of great baseball books"
$data[24,1] = "A new biography of Clayton Kershaw, a history of Chavez Ravine, a reevaluation of a Negro League star and more book recommendations for baseball fans."
$data[24,2] = "July 13, 2024"
$data[24,3] = "c61bf178-e69a-43c0-a1d6-6654cdf04d9e"
$data[24,4] = 0
$data[24,5] = $false

$data[25,0] = "Pricey camps. Family favors. Early dashes from work. How do parents survive summer?"
$data[25,1] = "The scramble for summer child care is no easy feat for working parents. Families must weave together summer camps, friendly favors and leeway at work — often at high cost."
$data[25,2] = "July 12, 2024"
$data[25,3] = "7bbe742a-17f9-448e-b699-1fe49461f450"
$data[25,4] = 0
$data[25,5] = $false

$data[26,0] = "Cal Poly Humboldt president to step down months after campus crackdowns"
$data[26,1] = "Cal Poly Humboldt University President Tom Jackson Jr. announced Thursday that he is stepping down after a five-year run that included national attention this spring over a campus crackdown on pro-Palestinian protests."
$data[26,2] = "July 12, 2024"
$data[26,3] = "93412ffc-c1a0-4a30-878d-44402d8941e1"
$data[26,4] = 0
$data[26,5] = $false

$data[27,0] = "Editorial: L.A. County shows L.A. city that reform doesn’t have to be hard"
$data[27,1] = "While the L.A. City Council dithered over the details and delayed changes that could affect their power, county leaders forged ahead with a comprehensive governance reform ballot proposal."
$data[27,2] = "July 11, 2024"
$data[27,3] = "2eb4fd54-e6de-4eca-894d-71445666306f"
$data[27,4] = 0
$data[27,5] = $false

$data[28,0] = "Supt. Carvalho moves ahead with troubled AI effort despite collapse of tech contractor"
$data[28,1] = "LAUSD unplugged chatbot after collapse of company that created it, but Carvalho says other features remain to help students. Most schools don’t yet have it."
$data[28,2] = "July 11, 2024"
$data[28,3] = "a9850ec3-2d59-4127-9728-29d66ab52ae7"
$data[28,4] = 0
$data[28,5] = $false

$data[29,0] = "Column: Investing through index funds is more popular than ever, so why is it becoming controversial?"
$data[29,1] = "More Americans own stocks than ever before, yet the passive index funds they favor are facing increasing partisan criticism. Here’s why"
$data[29,2] = "July 10, 2024"
$data[29,3] = "5bec0f87-6b06-496c-b27d-6d010db2c6cf"
$data[29,4] = 0
$data[29,5] = $false

$data[30,0] = "USC President Carol Folt’s contract is renewed, but university won’t say for how long"
$data[30,1] = "USC President Carol Folt, who has drawn praise for expanding student programs and boosting athletics but criticism for her handling of pro-Palestinian protests, has received an extension to her five-year contract by university trustees."
$data[30,2] = "July 09, 2024"
$data[30,3] = "b711a767-39e3-414b-9d3c-c4a12c7a35ee"
$data[30,4] = 0
$data[30,5] = $false

$data[31,0] = "Column: A Trump judge blocks another pro-worker Biden initiative, this one involving noncompete clauses"
$data[31,1] = "The Federal Trade Commission banned non-compete clauses, which block workers from moving to better jobs. A Trump-appointed judge has blocked it--of course"
$data[31,2] = "July 09, 2024"
$data[31,3] = "d87bdd3a-b65f-4c5b-9810-96ad0b179301"
$data[31,4] = 0
$data[31,5] = $false

$data[32,0] = "Meet the Californians serving in the first class of the American Climate Corps"
$data[32,1] = "The White House has sworn in more than 9,000 members of the American Climate Corps. In California, they’re managing wildfires, installing solar panels and more."
$data[32,2] = "July 09, 2024"
$data[32,3] = "b9b63e4a-ab90-42d6-b46c-9a29e4092645"
$data[32,4] = 0
$data[32,5] = $false

$data[33,0] = "Alberto Carvalho: Bold post-COVID-lockdown school leader"
$data[33,1] = "Alberto Carvalho brings a big resume and big style to the Los Angeles Unified School District. High hopes ride on him as students emerge from the COVID-19 pandemic."
$data[33,2] = "July 07, 2024"
$data[33,3] = "114eb2a7-9545-4bec-abb0-7875160bc445"
$data[33,4] = 0
$data[33,5] = $false

$data[34,0] = "Your guide to Proposition 2: Education bond"
$data[34,1] = "School bond supports say the money is direly needed to help fund repairs and upgrades at thousands of California public elementary, middle and high schools and community colleges."
$data[34,2] = "July 05, 2024"
$data[34,3] = "39fa7719-691c-45ef-85b1-39dee11e7311"
$data[34,4] = 0
$data[34,5] = $false

$data[35,0] = "Biden vows to keep running as signs point to rapidly eroding support on Capitol Hill"
$data[35,1] = "President Biden vows to keep running for reelection, rejecting pressure from within his Democratic Party to withdraw after his poor debate performance."
$data[35,2] = "July 03, 2024"
$data[35,3] = "28d786bc-56c1-4df2-ae7b-dd56617ae09b"
$data[35,4] = 0
$data[35,5] = $false

$data[36,0] = "LAUSD shelves its hyped AI chatbot to help students after collapse of firm that made it"
$data[36,1] = "LAUSD sidelines “Ed,” an AI chatbot, after a splashy kick-off featuring a company that has now tanked. District also is dealing with another data breach."
$data[36,2] = "July 03, 2024"
$data[36,3] = "e1375478-7b25-4bb3-9e21-39b3c0939df9"
$data[36,4] = 0
$data[36,5] = $false

$data[37,0] = "Opinion: Fentanyl could fuel another cycle of loss in L.A.’s Black communities. It doesn’t have to"
$data[37,1] = "Opioid overdoses are depriving hundreds of thousands of children of parents, with especially dire consequences for Black families in Southern California."
$data[37,2] = "July 02, 2024"
$data[37,3] = "dbdb6af0-cf99-44be-973b-73b9e63494c5"
$data[37,4] = 0
$data[37,5] = $false

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6))
$rng.Value = $data

Write-Host "Done: wrote $rowCount rows"